$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff shifts the Start Date / Fest Name / Fest Type / College columns
# (B:E) down by one data row (old row 2 -> row 3, ..., old row 10 -> row 11)
# and writes a brand-new entry into row 2. The old row 11 entry is
# discarded. Column A (the numeric index) is untouched by the diff.
# Apply as direct value writes, working from the bottom up so a row's
# original value is captured before it gets overwritten.

for ($r = 11; $r -ge 3; $r--) {
    $prev = $r - 1
    $ws.Range("B$r").Value2 = $ws.Range("B$prev").Value2
    $ws.Range("C$r").Value2 = $ws.Range("C$prev").Value2
    $ws.Range("D$r").Value2 = $ws.Range("D$prev").Value2
    $ws.Range("E$r").Value2 = $ws.Range("E$prev").Value2
}

$ws.Range("B2").Value2 = "28 Feb 2023"
$ws.Range("C2").Value2 = " Smart Bengal Hackathon 2023 "
$ws.Range("D2").Value2 = "Technical, Hackathon"
$ws.Range("E2").Value2 = "RCC Institute of Information Technology"
